$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")

# --- Settings sheet edits (row 10-16) ---
# NOTE: cell writes are deliberately ordered to reproduce the exact shared-string
# insertion order seen in the target workbook (the author apparently filled in
# the "UpcomingBirthday_*" rows (14-16) before the "BirthdayList_Index_*" rows
# (10-12), and within that second group wrote the EID row before the Birthday row).

# Row 14: rename existing entry (reuses description text already in the table)
$settings.Range("A14").Value = "UpcomingBirthday_Range"
$settings.Range("B14").Value = 14
$settings.Range("C14").Value = "Value in number of days to determine upcoming birthdays within range (DateToday - X days)"

# Row 15 (new)
$settings.Range("A15").Value = "UpcomingBirthday_OutputFolder"
$settings.Range("B15").Value = "C:\Users\{{UserProfile}}\Documents\UiPath\Birthday Greeting Postcard\Upcoming birthdays"
$settings.Range("C15").Value = "Filepath of the upcoming birthday celebrant's team mates"

# Row 16 (new)
$settings.Range("A16").Value = "UpcomingBirthday_IncludeTeam"
$settings.Range("B16").Value = "Management"
$settings.Range("C16").Value = 'Included team when consolidating the celebrant''s team members. Only accepts 1 team. Write "n/a" if you do not wish to include other team.'

# Row 11 (new) - written before row 10 to match original authoring order
$settings.Range("A11").Value = "BirthdayList_Index_EID"
$settings.Range("B11").Value = 4

# Row 10 (overwrite old "UpcomingBirthdayRange" entry)
$settings.Range("A10").Value = "BirthdayList_Index_Birthday"
$settings.Range("B10").Value = 3

# Row 12 (new)
$settings.Range("A12").Value = "BirthdayList_Index_Team"
$settings.Range("B12").Value = 7

# Descriptions (col C) written after all the names/values above
$settings.Range("C10").Value = 'Column index of "Birthday" in BirthdayList datatable'
$settings.Range("C11").Value = 'Column index of "Birthday" in BirthdayList datatable'
$settings.Range("C12").Value = 'Column index of "Team" in BirthdayList datatable'

# Extend used range / formatted rows at bottom (999-1002) to match new dimension
$settings.Rows("999:1002").RowHeight = 14.25

# --- Sheet view / selection updates ---
$settings.Range("B13").Select()

# Make Settings the active/selected sheet, Constants no longer the tab-selected one
$settings.Activate()

$wb.Application.ActiveWindow.WindowState = $wb.Application.ActiveWindow.WindowState
